$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting of the existing "sum" header cell (G1) onto the new
# "Save" header cell (H1) so it reuses the same style index as the rest
# of the header row.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Set the new header text and the value for row 2.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
